$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 83; $row++) {
    $cell = $ws.Range("C$row")
    $cell.Value2 = $cell.Value2 + 1
}
